$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6021.32
$ws.Range("J64").Value = 7699.125
$ws.Range("L64").Value = 7699.125
$ws.Range("N64").Value = -8195.125

$ws.Range("H67").Value = 6021.32
$ws.Range("J67").Value = 7699.125
$ws.Range("L67").Value = 7699.125
$ws.Range("N67").Value = -9415.125

$ws.Range("H103").Value = 2065.5715
$ws.Range("J103").Value = 1593
$ws.Range("L103").Value = 4779
$ws.Range("N103").Value = -5951

$ws.Range("H125").Value = 6369.3335
$ws.Range("I125").Value = 6722
$ws.Range("K125").Value = 60498
$ws.Range("M125").Value = -58038

$ws.Range("H137").Value = 606909.4399999999
$ws.Range("J137").Value = 1615332.5
$ws.Range("L137").Value = 4845997.5
$ws.Range("N137").Value = -4851097.5


# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7300.635
$ws.Range("I32").Value = 2465.1516
$ws.Range("K32").Value = 2465.1516
$ws.Range("M32").Value = -2178.1516

$ws.Range("H63").Value = 1818.091
$ws.Range("I63").Value = 1874.9
$ws.Range("J63").Value = 1250
$ws.Range("K63").Value = 1874.9
$ws.Range("L63").Value = 1250
$ws.Range("M63").Value = -1188.9
$ws.Range("N63").Value = -2622

$ws.Range("H66").Value = 1818.091
$ws.Range("I66").Value = 1874.9
$ws.Range("J66").Value = 1250
$ws.Range("K66").Value = 9374.5
$ws.Range("L66").Value = 6250
$ws.Range("M66").Value = -5942.5
$ws.Range("N66").Value = -13114


# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2642428.2
$ws.Range("I99").Value = 126231.625
$ws.Range("K99").Value = 126231.625
$ws.Range("M99").Value = -124733.625


# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3744.4
$ws.Range("I31").Value = 3080.7693
$ws.Range("J31").Value = 4463.3335
$ws.Range("K31").Value = 3080.7693
$ws.Range("L31").Value = 4463.3335
$ws.Range("M31").Value = -2785.7693
$ws.Range("N31").Value = -5053.3335

$ws.Range("H33").Value = 6575.8
$ws.Range("I33").Value = 1440.5
$ws.Range("J33").Value = 9999.333000000001
$ws.Range("K33").Value = 1440.5
$ws.Range("L33").Value = 9999.333000000001
$ws.Range("M33").Value = -1061.5
$ws.Range("N33").Value = -10757.333

$ws.Range("H34").Value = 3744.4
$ws.Range("I34").Value = 3080.7693
$ws.Range("J34").Value = 4463.3335
$ws.Range("K34").Value = 3080.7693
$ws.Range("L34").Value = 4463.3335
$ws.Range("M34").Value = -2878.7693
$ws.Range("N34").Value = -4867.3335

$ws.Range("H132").Value = 2089765.2
$ws.Range("I132").Value = 2068223.5
$ws.Range("J132").Value = 2168751.2
$ws.Range("K132").Value = 6204670.5
$ws.Range("L132").Value = 6506253.600000001
$ws.Range("M132").Value = -6202140.5
$ws.Range("N132").Value = -6511313.600000001


# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 792.64
$ws.Range("I5").Value = 644.2857
$ws.Range("J5").Value = 981.4545000000001
$ws.Range("K5").Value = 1932.8571
$ws.Range("L5").Value = 2944.3635
$ws.Range("M5").Value = -1820.8571
$ws.Range("N5").Value = -3168.3635

$ws.Range("H113").Value = 2757565.5
$ws.Range("J113").Value = 3545014.5
$ws.Range("L113").Value = 10635043.5
$ws.Range("N113").Value = -10639383.5

$ws.Range("H121").Value = 2376.647
$ws.Range("I121").Value = 983.4
$ws.Range("J121").Value = 2957.1667
$ws.Range("K121").Value = 2950.2
$ws.Range("L121").Value = 8871.500100000001
$ws.Range("M121").Value = -1640.2
$ws.Range("N121").Value = -11491.5001

$ws.Range("H135").Value = 792.64
$ws.Range("I135").Value = 644.2857
$ws.Range("J135").Value = 981.4545000000001
$ws.Range("K135").Value = 5798.571300000001
$ws.Range("L135").Value = 8833.0905
$ws.Range("M135").Value = -3263.571300000001
$ws.Range("N135").Value = -13903.0905


# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 24725
$ws.Range("I57").Value = 9450
$ws.Range("K57").Value = 9450
$ws.Range("M57").Value = -8630

$ws.Range("H62").Value = 27563.75

$ws.Range("H65").Value = 27563.75

$ws.Range("H70").Value = 6354.647
$ws.Range("I70").Value = 6493.5
$ws.Range("J70").Value = 6231.222
$ws.Range("K70").Value = 6493.5
$ws.Range("L70").Value = 6231.222
$ws.Range("M70").Value = -6223.5
$ws.Range("N70").Value = -6771.222

$ws.Range("H73").Value = 6354.647
$ws.Range("I73").Value = 6493.5
$ws.Range("J73").Value = 6231.222
$ws.Range("K73").Value = 6493.5
$ws.Range("L73").Value = 6231.222
$ws.Range("M73").Value = -5557.5
$ws.Range("N73").Value = -8103.222

$ws.Range("H132").Value = 8063.647
$ws.Range("I132").Value = 7167
$ws.Range("J132").Value = 8437.25
$ws.Range("K132").Value = 21501
$ws.Range("L132").Value = 25311.75
$ws.Range("M132").Value = -18971
$ws.Range("N132").Value = -30371.75


# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 200.8

$ws.Range("H5").Value = 6743.1665
$ws.Range("I5").Value = 4986.3335
$ws.Range("J5").Value = 8500
$ws.Range("K5").Value = 4986.3335
$ws.Range("L5").Value = 8500
$ws.Range("M5").Value = -4873.3335
$ws.Range("N5").Value = -8726

$ws.Range("H7").Value = 5299.8335
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 5299.8335
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 5299.8335
$ws.Range("N7").Value = -5523.8335
$ws.Range("M7").ClearContents()

$ws.Range("H56").Value = 8450
$ws.Range("J56").Value = 17000
$ws.Range("L56").Value = 17000
$ws.Range("N56").Value = -18382

$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H126").Value = 5299.8335
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5299.8335
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15899.5005
$ws.Range("N126").Value = -20839.5005
$ws.Range("M126").ClearContents()


# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 19533.75
$ws.Range("I51").Value = 19533.75
$ws.Range("K51").Value = 19533.75
$ws.Range("M51").Value = -19023.75

$ws.Range("H58").Value = 56194
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 56194
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 56194
$ws.Range("N58").Value = -56810
$ws.Range("M58").ClearContents()

$ws.Range("H100").Value = 3969115
$ws.Range("I100").Value = 5495426
$ws.Range("J100").Value = 706.4
$ws.Range("K100").Value = 10990852
$ws.Range("L100").Value = 1412.8
$ws.Range("M100").Value = -10990311
$ws.Range("N100").Value = -2494.8

$ws.Range("H113").Value = 2260.0908
$ws.Range("I113").Value = 2830.8333
$ws.Range("K113").Value = 8492.499899999999
$ws.Range("M113").Value = -6322.499899999999

$ws.Range("H122").Value = 1919.7693
$ws.Range("I122").Value = 1961.6
$ws.Range("J122").Value = 1893.625
$ws.Range("K122").Value = 5884.799999999999
$ws.Range("L122").Value = 5680.875
$ws.Range("M122").Value = -3434.799999999999
$ws.Range("N122").Value = -10580.875

